# Email Farma Conde config workbook update:
#  - expand the "emailCc" list (cell C2) with the new Farma Conde recipients
#  - remove the stray/blank formatted cell that had been left in C3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CC e-mail list held in C2 (semicolon separated, no spaces,
# with the new Farma Conde recipients appended).
$ws.Range("C2").Value = "a.chagas@senff.com.br;l.ruiz@senff.com.br;samir.nadir@grupofarmaconde.com.br;danilo.fernandes@grupofarmaconde.com.br;valeska.amorim@grupofarmaconde.com.br"

# Remove the extra (empty, bold-styled) row 3 / cell C3 that was left over
# in the sheet, shrinking the used range back down to A1:F2.
$ws.Rows.Item(3).Delete()

# Reflect the resulting selection in the sheet view (the cursor had moved
# down past the removed row).
$ws.Range("C9").Select()
